$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Fuscoporia viticola" (Vedticka), "Phellinidium
# ferrugineofuscum" (Ullticka), and "Polygonatum verticillatum" (Kransrams)
# were cyclically rotated: row2 <- old row3, row3 <- old row4, row4 <- old row2
# (columns A, B, D, E, F, G, H; and Q/R go along with the record they describe).

$row2 = @{ A = 111463862; B = 89369;  D = "LC"; E = 5447;   F = "Vedticka";  G = "Fuscoporia viticola";              H = "(Schwein.) Murrill";               Q = 554109.1038748255;  R = 7007938.027731327 }
$row3 = @{ A = 111463857; B = 89405;  D = "NT"; E = 1202;   F = "Ullticka";  G = "Phellinidium ferrugineofuscum";    H = "(P.Karst.) Fiasson & Niemelä";     Q = 554109.1038748255;  R = 7007938.027731327 }
$row4 = @{ A = 111463670; B = 96674;  D = "LC"; E = 219880; F = "Kransrams"; G = "Polygonatum verticillatum";        H = "(L.) All.";                        Q = 554151.0634843309;  R = 7007942.793868498 }

$new2 = $row3
$new3 = $row4
$new4 = $row2

function Set-RowValues($rowIndex, $data) {
    $ws.Range("A$rowIndex").Value = $data.A
    $ws.Range("B$rowIndex").Value = $data.B
    $ws.Range("D$rowIndex").Value = $data.D
    $ws.Range("E$rowIndex").Value = $data.E
    $ws.Range("F$rowIndex").Value = $data.F
    $ws.Range("G$rowIndex").Value = $data.G
    $ws.Range("H$rowIndex").Value = $data.H
    $ws.Range("Q$rowIndex").Value = $data.Q
    $ws.Range("R$rowIndex").Value = $data.R
}

Set-RowValues 2 $new2
Set-RowValues 3 $new3
Set-RowValues 4 $new4
